$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: A2 1 -> 0, B2 1 -> 3
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 3

# Remove rows 3 and 4 entirely (they held data that is no longer needed)
$ws.Rows("3:4").Delete()
